$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-06-10 Tuesday" "2025-06-11 Wednesday"

Replace-Text "93×91=8463" "35×94=3290"
Replace-Text "46×39=1794" "98×88=8624"
Replace-Text "31×50=1550" "61×91=5551"
Replace-Text "88×69=6072" "92×26=2392"
Replace-Text "41×26=1066" "59×49=2891"

Replace-Text "27×18=486" "19×26=494"
Replace-Text "25×54=1350" "60×40=2400"
Replace-Text "25×43=1075" "35×54=1890"
Replace-Text "18×84=1512" "80×90=7200"
Replace-Text "96×30=2880" "39×85=3315"

Replace-Text "89×94=8366" "62×39=2418"
Replace-Text "80×32=2560" "31×16=496"
Replace-Text "17×19=323" "98×57=5586"
Replace-Text "55×57=3135" "62×70=4340"
Replace-Text "58×52=3016" "67×40=2680"

Replace-Text "14×79=1106" "93×94=8742"
Replace-Text "46×42=1932" "17×22=374"
Replace-Text "25×73=1825" "23×34=782"
Replace-Text "11×23=253" "34×19=646"
Replace-Text "26×82=2132" "54×30=1620"

Replace-Text "62×24=1488" "76×86=6536"
Replace-Text "12×66=792" "43×15=645"
Replace-Text "96×80=7680" "86×79=6794"
Replace-Text "76×98=7448" "36×82=2952"
Replace-Text "86×99=8514" "74×88=6512"
